$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.378.97"
$ws.Range("E2").Value = "  -7.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.895.37"
$ws.Range("E3").Value = "  -10.44%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "473.09"
$ws.Range("E5").Value = "  -12.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.44"
$ws.Range("E6").Value = "  -7.47%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.897.79"
$ws.Range("E8").Value = "  -10.35%  "

$ws.Range("E9").Value = "  -12.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  -12.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0972"
$ws.Range("E11").Value = "  -15.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.332"
$ws.Range("E12").Value = "  -15.89%  "

$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.396.66"
$ws.Range("E14").Value = "  -10.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.99"
$ws.Range("E15").Value = "  -11.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.502.01"
$ws.Range("E16").Value = "  -7.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.898.49"
$ws.Range("E17").Value = "  -10.42%  "

$ws.Range("E18").Value = "  -15.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.17"
$ws.Range("E19").Value = "  -12.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -13.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.16"
$ws.Range("E21").Value = "  -13.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "308.57"
$ws.Range("E22").Value = "  -14.78%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.449"
$ws.Range("E24").Value = "  -13.89%  "

$ws.Range("E25").Value = "  -16.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -9.25%  "

$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0829"
$ws.Range("E29").Value = "  -15.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.23"
$ws.Range("E30").Value = "  -12.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  -6.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  -13.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.22"
$ws.Range("E33").Value = "  -12.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.62"
$ws.Range("E34").Value = "  -16.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "140.74"
$ws.Range("E35").Value = "  -13.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.24"
$ws.Range("E36").Value = "  -14.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.47"
$ws.Range("E37").Value = "  -15.06%  "

$ws.Range("E38").Value = "  -15.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.12"
$ws.Range("E39").Value = "  -12.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0623"
$ws.Range("E40").Value = "  -12.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.925.60"
$ws.Range("E41").Value = "  -10.42%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.76"
$ws.Range("E43").Value = "  -15.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.971"
$ws.Range("E44").Value = "  -12.24%  "

$ws.Range("E45").Value = "  -16.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.45"
$ws.Range("E46").Value = "  -14.68%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.066.26"
$ws.Range("E47").Value = "  -10.15%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("E48").Value = "  -13.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.19"
$ws.Range("E49").Value = "  -13.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.32"
$ws.Range("E50").Value = "  -15.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0213"
$ws.Range("E51").Value = "  -12.16%  "
